# Refresh the cryptocurrency price/volume snapshot (and one reshuffled
# ranking swap between Aave and BabyDogeCoin) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume (and, for two rows, Coin/Link) cells are stored as
# plain text in the source data (e.g. "28.907.13", "  -0.43%  "). Some
# of the new values parse as ordinary numbers (e.g. "8.520"), which would
# otherwise be auto-coerced to a numeric cell and lose formatting such as
# trailing zeros. Prefixing the value with a literal quote-prefix keeps it
# as text, matching the original cell type; resetting the style back to
# "Normal" afterwards avoids leaving a visible text-number-format behind.
function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "29.023.39"
Set-TextValue "E2" "  -0.03%  "
Set-TextValue "D3" "1.829.53"
Set-TextValue "E3" "  -0.16%  "
Set-TextValue "D4" "0.9979"
Set-TextValue "E4" "  -0.19%  "
Set-TextValue "D5" "243.76"
Set-TextValue "E5" "  +0.94%  "
Set-TextValue "D6" "0.6316"
Set-TextValue "E6" "  +0.70%  "
Set-TextValue "D7" "0.9991"
Set-TextValue "E7" "  -0.16%  "
Set-TextValue "D8" "0.07494"
Set-TextValue "E8" "  -1.03%  "
Set-TextValue "D9" "0.2943"
Set-TextValue "E9" "  +0.87%  "
Set-TextValue "D10" "23.09"
Set-TextValue "E10" "  +1.20%  "
Set-TextValue "D11" "0.07706"
Set-TextValue "E11" "  +0.78%  "
Set-TextValue "D12" "1.825.29"
Set-TextValue "E12" "  -0.11%  "
Set-TextValue "D13" "4.998"
Set-TextValue "E13" "  +0.91%  "
Set-TextValue "D14" "0.6689"
Set-TextValue "E14" "  +0.48%  "
Set-TextValue "D15" "83.15"
Set-TextValue "E15" "  +0.90%  "
Set-TextValue "D16" "0.000009812"
Set-TextValue "E16" "  +3.91%  "
Set-TextValue "D17" "6.039"
Set-TextValue "E17" "  +0.84%  "
Set-TextValue "D18" "29.012.56"
Set-TextValue "E18" "  -0.01%  "
Set-TextValue "D19" "12.58"
Set-TextValue "E19" "  +2.06%  "
Set-TextValue "D20" "226.22"
Set-TextValue "E20" "  +0.58%  "
Set-TextValue "D21" "0.9988"
Set-TextValue "E21" "  -0.10%  "
Set-TextValue "D22" "7.135"
Set-TextValue "E22" "  -0.90%  "
Set-TextValue "D23" "0.9994"
Set-TextValue "E23" "  -0.20%  "
Set-TextValue "D24" "160.71"
Set-TextValue "E24" "  +0.43%  "
Set-TextValue "D25" "0.1414"
Set-TextValue "E25" "  +3.64%  "
Set-TextValue "D26" "8.520"
Set-TextValue "E26" "  +1.26%  "
Set-TextValue "D27" "17.95"
Set-TextValue "E27" "  +0.69%  "
Set-TextValue "D28" "1.501"
Set-TextValue "E28" "  +0.48%  "
Set-TextValue "D29" "4.130"
Set-TextValue "E29" "  +1.67%  "
Set-TextValue "D30" "4.057"
Set-TextValue "E30" "  +0.67%  "
Set-TextValue "D31" "0.05485"
Set-TextValue "E31" "  +5.43%  "
Set-TextValue "D32" "1.201"
Set-TextValue "E32" "  +0.27%  "
Set-TextValue "D33" "1.859"
Set-TextValue "E33" "  +0.67%  "
Set-TextValue "D34" "0.7453"
Set-TextValue "E34" "  +1.80%  "
Set-TextValue "D35" "1.138"
Set-TextValue "E35" "  -1.14%  "
Set-TextValue "D36" "2.611"
Set-TextValue "E36" "  +1.23%  "
Set-TextValue "D37" "1.243.29"
Set-TextValue "E37" "  -1.96%  "
Set-TextValue "D38" "2.749"
Set-TextValue "E38" "  -0.35%  "
Set-TextValue "D39" "0.01782"
Set-TextValue "E39" "  -0.27%  "
Set-TextValue "D40" "6.712"
Set-TextValue "E40" "  +2.90%  "
Set-TextValue "D41" "0.9019"
Set-TextValue "E41" "  +1.18%  "
Set-TextValue "D42" "0.9991"
Set-TextValue "E42" "  -0.13%  "
Set-TextValue "D43" "101.57"
Set-TextValue "E43" "  -0.08%  "
Set-TextValue "D44" "1.968.18"
Set-TextValue "E44" "  -0.48%  "
Set-TextValue "B45" "BabyDogeCoin"
Set-TextValue "C45" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D45" "0.00000000124"
Set-TextValue "E45" "  +2.96%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "65.21"
Set-TextValue "E46" "  +0.62%  "
Set-TextValue "D47" "0.5078"
Set-TextValue "E47" "  -0.64%  "
Set-TextValue "D48" "0.4052"
Set-TextValue "E48" "  +1.89%  "
Set-TextValue "D49" "0.07455"
Set-TextValue "E49" "  +4.80%  "
Set-TextValue "D50" "8.955"
Set-TextValue "E50" "  +1.15%  "
Set-TextValue "D51" "1.666"
Set-TextValue "E51" "  +1.86%  "
